# Apply the updates described by the diff:
# 1. Update the shared "IndicatorQuantiles.R, Git Commit ID: ..." text (column AJ, rows 2-80)
#    from the old commit hash to the new commit hash.
# 2. Update the "pid" values in column AH (rows 2-80) from 25596 to 25080.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldCommitText = "IndicatorQuantiles.R, Git Commit ID: 0e4152332be22faf035a2e2fc83ad2cca4c8a7fc"
$newCommitText = "IndicatorQuantiles.R, Git Commit ID: 2e3ff9a54734c37c56b32bb788c6f054c2509b6b"

$lastRow = 80

for ($row = 2; $row -le $lastRow; $row++) {
    # Column AH = 34 ("pid") : 25596 -> 25080
    $ahCell = $ws.Cells.Item($row, 34)
    if ($ahCell.Value() -eq 25596) {
        $ahCell.Value = 25080
    }

    # Column AJ = 36 ("ScriptLatestRunVersion") : update commit id text
    $ajCell = $ws.Cells.Item($row, 36)
    if ($ajCell.Value() -eq $oldCommitText) {
        $ajCell.Value = $newCommitText
    }
}
